$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I ("Antal") holds numeric-looking values stored as text -
# force text format on just the edited cells so COM doesn't reinterpret
# them as numbers (and we don't touch untouched rows like I5).
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I6").NumberFormat = "@"

# New values for row 2 (previously row 6's values for these columns)
$ws.Range("A2").Value = 111416526
$ws.Range("I2").Value = "1"
$ws.Range("Q2").Value = 359094.3997885482
$ws.Range("R2").Value = 6393206.775113393
$ws.Range("AC2").Value = "även ca 30 bladrosetter"

# New values for row 3 (previously row 4's values for these columns)
$ws.Range("A3").Value = 111416525
$ws.Range("I3").Value = "4"
$ws.Range("Q3").Value = 359095.1406046218
$ws.Range("R3").Value = 6393212.639220579
$ws.Range("AC3").Value = "även bladrosetter på 1 kvm"

# New values for row 4 (previously row 3's values for these columns)
$ws.Range("A4").Value = 111416521
$ws.Range("I4").Value = "1"
$ws.Range("Q4").Value = 359101.3469427949
$ws.Range("R4").Value = 6393205.997596246
$ws.Range("AC4").Value = "även ca 30 bladrosetter"

# New values for row 6 (previously row 2's values for these columns)
$ws.Range("A6").Value = 111416523
$ws.Range("I6").Value = "7"
$ws.Range("Q6").Value = 359100.0376043977
$ws.Range("R6").Value = 6393214.610374114
$ws.Range("AC6").Value = "även bladrosetter på 1 kvm"
